# Generate Report for Handoff
# - Priority changes from "low" to "ht" for rows 4-7 (215d6b56, 2b9ae4b8, 51568b8e, a8e133f9)
#   on both the zh-cn and de-de sheets.
# - The zh-cn "Latest Handoff Datetime" for those same rows updates from
#   2016-09-03 14:34:02 to 2016-09-03 14:34:24.
# - The shared "Latest HO Xliff Generate Date" / de-de "Latest Handoff Datetime"
#   value for the 215d6b56 file updates from 2016-09-03 14:34:11 to 2016-09-03 14:34:28
#   (this value is shared across the Overview sheet and the de-de sheet).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Priority column (E) on rows 4-7: "low" -> "ht"
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# zh-cn Latest Handoff Datetime (H) rows 4-7: 2016-09-03 14:34:02 -> 2016-09-03 14:34:24
foreach ($row in 4..7) {
    $zhcn.Range("H$row").Value = "2016-09-03 14:34:24"
}

# Shared timestamp 2016-09-03 14:34:11 -> 2016-09-03 14:34:28
# This backs the Overview "Latest HO Xliff Generate Date" column (G) for every row
# as well as the de-de "Latest Handoff Datetime" column (H) for rows 4-7.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-09-03 14:34:28"
    $dede.Range("H$row").Value = "2016-09-03 14:34:28"
}
